# Generate Report for Handoff
# Adds a new tracked file (a98e7d6d-...) as row 3 on the Overview, zh-cn
# and de-de sheets, reflecting that it is "Ready for handoff".

$wb = $excel.ActiveWorkbook

$newFile        = "a98e7d6d-d44b-4c14-b428-0cd1edd29409ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFilePath    = "e2e\a98e7d6d-d44b-4c14-b428-0cd1edd29409ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$statusText     = "Ready for handoff"
$overviewDate   = "2016-08-21 12:35:55"
$zhDate         = "2016-08-21 12:35:51"
$deDate         = "2016-08-21 12:35:55"
$zhXliff        = "a98e7d6d-d44b-4c14-b428-0cd1edd29409oooooooooooooooooooooooooooooooooooooooo.d84d4889c28b2a42b0e83271c5bd23e9c6f64749.zh-cn.xlf"
$deXliff        = "a98e7d6d-d44b-4c14-b428-0cd1edd29409oooooooooooooooooooooooooooooooooooooooo.d84d4889c28b2a42b0e83271c5bd23e9c6f64749.de-de.xlf"
$epoch          = "0001-01-01 00:00:00"

$commitHash     = "fae63f150e1bf098c6ba2354a895c1e5c1ec5f2b"
$newFileUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $newFile

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = $statusText
$wsOverview.Cells.Item(3, 6).Value = $statusText
$wsOverview.Cells.Item(3, 7).Value = $overviewDate
$wsOverview.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 1), $newFileUrl, "", "", $newFile)
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), $newFileUrl, "", "", $newFilePath)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = $statusText
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "False"
$wsZh.Cells.Item(3, 8).Value = $zhDate
$wsZh.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 11).Value = $epoch
$wsZh.Cells.Item(3, 13).Value = "True"
$wsZh.Cells.Item(3, 15).Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 1), $newFileUrl, "", "", $newFile)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 7), $newFileUrl, "", "", $zhXliff)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = $statusText
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "False"
$wsDe.Cells.Item(3, 8).Value = $deDate
$wsDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 11).Value = $epoch
$wsDe.Cells.Item(3, 13).Value = "True"
$wsDe.Cells.Item(3, 15).Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 1), $newFileUrl, "", "", $newFile)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 7), $newFileUrl, "", "", $deXliff)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
